# Update the chart number entry in B2 to include the drive-number suffix
# (e.g. "64275" -> "64275_1"), per the diff: sheet1.xml cell B2 now
# references a new shared string "64275_1" instead of the inline string
# "64275".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = "64275_1"
